# Updated cryptos list (price + 1h volume-change columns) as captured by the
# GitHub Actions scraper run. Price-looking values that Excel would otherwise
# auto-convert to numbers are entered with a leading apostrophe so they stay
# stored as text, matching the original column formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.109.35"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "3.274.71"

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'575.84"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").Value = "'179.26"
$ws.Range("E6").Value = "  -3.02%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +3.20%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "3.848.02"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("E13").Value = "  -3.88%  "

$ws.Range("D14").Value = "66.153.22"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").Value = "'26.36"
$ws.Range("E15").Value = "  -3.43%  "

$ws.Range("D16").Value = "'0.0000162"
$ws.Range("E16").Value = "  -2.56%  "

$ws.Range("D17").Value = "3.235.59"
$ws.Range("E17").Value = "  -2.51%  "

$ws.Range("D18").Value = "'430.09"
$ws.Range("E18").Value = "  -2.89%  "

$ws.Range("D19").Value = "'5.53"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("D20").Value = "'13.07"
$ws.Range("E20").Value = "  -3.61%  "

$ws.Range("D21").Value = "'7.37"
$ws.Range("E21").Value = "  -4.34%  "

$ws.Range("D22").Value = "'71.86"
$ws.Range("E22").Value = "  -3.09%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "3.417.82"
$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").Value = "'0.503"
$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("E26").Value = "  +3.58%  "

$ws.Range("D27").Value = "'0.0000112"
$ws.Range("E27").Value = "  -5.19%  "

$ws.Range("D28").Value = "'8.81"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -1.98%  "

$ws.Range("D31").Value = "'22.21"
$ws.Range("E31").Value = "  -3.10%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").Value = "'5.14"
$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("D34").Value = "'6.55"
$ws.Range("E34").Value = "  -3.82%  "

$ws.Range("E35").Value = "  -4.41%  "

$ws.Range("D36").Value = "'157.64"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("D37").Value = "'1.42"
$ws.Range("E37").Value = "  -5.69%  "

$ws.Range("D38").Value = "'26.39"
$ws.Range("E38").Value = "  -3.10%  "

$ws.Range("D39").Value = "'1.78"
$ws.Range("E39").Value = "  -3.39%  "

$ws.Range("D40").Value = "2.758.64"
$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("D41").Value = "'0.774"
$ws.Range("E41").Value = "  -1.82%  "

$ws.Range("D42").Value = "'4.28"
$ws.Range("E42").Value = "  -4.30%  "

$ws.Range("D43").Value = "'40.25"

$ws.Range("D44").Value = "'6.02"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").Value = "'0.0655"

$ws.Range("E46").Value = "  -3.46%  "

$ws.Range("D47").Value = "'320.41"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").Value = "'23.10"
$ws.Range("E48").Value = "  -6.20%  "

$ws.Range("D49").Value = "'0.0265"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("E51").Value = "  +0.01%  "
